$d = $word.ActiveDocument

# --- Edit 1: First paragraph gets two trailing spaces, then a red-colored
#             parenthetical note appended as three separate runs. ---
$p1 = $d.Paragraphs.Item(1)
$rng = $p1.Range
$rng.End = $rng.End - 1
$rng.Collapse(0)
$rng.InsertAfter("  ")
$rng.Collapse(0)

$rng2 = $d.Range($rng.End, $rng.End)
$rng2.InsertAfter("(This is a change – Ve")
$rng2.Font.Color = 255
$rng2.Collapse(0)

$rng3 = $d.Range($rng2.End, $rng2.End)
$rng3.InsertAfter("rsion for main branch")
$rng3.Font.Color = 255
$rng3.Collapse(0)

$rng4 = $d.Range($rng3.End, $rng3.End)
$rng4.InsertAfter(")")
$rng4.Font.Color = 255
$rng4.Collapse(0)

# --- Edit 2: Remove the final paragraph ("ank God almighty, we are free
#             at last."), which trails the "Shall be lifted-nevermore!"
#             paragraph. ---
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$lastPara.Range.Delete()
